# PersData_Excel.xlsx update:
#  - rename sheets: Unit1 -> Year1, Unit2 -> Year2
#  - move the selection/active tab from Unit1!F8 to Unit2!D22
#    (Year2 becomes the activated / tab-selected sheet)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Year1"
$ws2.Name = "Year2"

# Year1 keeps its own selection (E7) but is no longer the visible/active tab.
$ws1.Activate()
$ws1.Range("E7").Select() | Out-Null

# Year2 becomes the active sheet, with the selection moved to D22.
$ws2.Activate()
$ws2.Range("D22").Select() | Out-Null
